$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1522.1111
$ws.Range("J17").Value = 1522.1111
$ws.Range("L17").Value = 4566.3333
$ws.Range("N17").Value = -4902.3333

$ws.Range("H19").Value = 2206
$ws.Range("I19").Value = 1125.1111
$ws.Range("K19").Value = 1125.1111
$ws.Range("M19").Value = -950.1111000000001

$ws.Range("H74").Value = 5225
$ws.Range("I74").Value = 5471.4287
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 5471.4287
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -4535.4287
$ws.Range("N74").Value = -5372

$ws.Range("H77").Value = 5225
$ws.Range("I77").Value = 5471.4287
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 27357.1435
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -22677.1435
$ws.Range("N77").Value = -26860

$ws.Range("H88").Value = 899711.25
$ws.Range("I88").Value = 4451.5
$ws.Range("J88").Value = 1155499.8
$ws.Range("K88").Value = 4451.5
$ws.Range("L88").Value = 1155499.8
$ws.Range("M88").Value = -4045.5
$ws.Range("N88").Value = -1156311.8

$ws.Range("H91").Value = 899711.25
$ws.Range("I91").Value = 4451.5
$ws.Range("J91").Value = 1155499.8
$ws.Range("K91").Value = 4451.5
$ws.Range("L91").Value = 1155499.8
$ws.Range("M91").Value = -3047.5
$ws.Range("N91").Value = -1158307.8

$ws.Range("H99").Value = 540.2
$ws.Range("I99").Value = 475.25
$ws.Range("J99").Value = 800
$ws.Range("K99").Value = 1425.75
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = 72.25
$ws.Range("N99").Value = -5396

$ws.Range("H132").Value = 2917.318
$ws.Range("I132").Value = 2779.805
$ws.Range("J132").Value = 4796.6665
$ws.Range("K132").Value = 8339.414999999999
$ws.Range("L132").Value = 14389.9995
$ws.Range("M132").Value = -5809.414999999999
$ws.Range("N132").Value = -19449.9995

$ws.Range("H138").Value = 3323.694
$ws.Range("I138").Value = 1883
$ws.Range("J138").Value = 3376.4023
$ws.Range("K138").Value = 5649
$ws.Range("L138").Value = 10129.2069
$ws.Range("M138").Value = -509
$ws.Range("N138").Value = -20409.2069

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 901.4286
$ws.Range("I4").Value = 562.6
$ws.Range("J4").Value = 1748.5
$ws.Range("K4").Value = 562.6
$ws.Range("L4").Value = 1748.5
$ws.Range("M4").Value = -446.6
$ws.Range("N4").Value = -1980.5

$ws.Range("H5").Value = 252.20589
$ws.Range("J5").Value = 250.72728
$ws.Range("L5").Value = 250.72728
$ws.Range("N5").Value = -474.72728

$ws.Range("H32").Value = 13161038
$ws.Range("I32").Value = 14287755
$ws.Range("K32").Value = 14287755
$ws.Range("M32").Value = -14287468

$ws.Range("H61").Value = 15663898
$ws.Range("I61").Value = 21742620
$ws.Range("K61").Value = 21742620
$ws.Range("M61").Value = -21742408

$ws.Range("H97").Value = 1751.591
$ws.Range("I97").Value = 1815.9524
$ws.Range("K97").Value = 1815.9524
$ws.Range("M97").Value = -1319.9524

$ws.Range("H122").Value = 3937
$ws.Range("I122").Value = 3332.8333
$ws.Range("J122").Value = 5749.5
$ws.Range("K122").Value = 9998.499899999999
$ws.Range("L122").Value = 17248.5
$ws.Range("M122").Value = -7548.499899999999
$ws.Range("N122").Value = -22148.5

$ws.Range("H132").Value = 2864.4348
$ws.Range("I132").Value = 1489.6757
$ws.Range("K132").Value = 4469.0271
$ws.Range("M132").Value = -1939.0271

$ws.Range("H136").Value = 15663898
$ws.Range("I136").Value = 21742620
$ws.Range("K136").Value = 65227860
$ws.Range("M136").Value = -65225310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 252.20589
$ws.Range("J4").Value = 250.72728
$ws.Range("L4").Value = 250.72728
$ws.Range("N4").Value = -480.72728

$ws.Range("H96").Value = 34641.934
$ws.Range("I96").Value = 7687.273
$ws.Range("J96").Value = 108767.25
$ws.Range("K96").Value = 7687.273
$ws.Range("L96").Value = 108767.25
$ws.Range("M96").Value = -4941.273
$ws.Range("N96").Value = -114259.25

$ws.Range("H97").Value = 25194.4
$ws.Range("J97").Value = 48989.5
$ws.Range("L97").Value = 48989.5
$ws.Range("N97").Value = -50971.5

$ws.Range("H134").Value = 418356.6
$ws.Range("I134").Value = 1454.5238
$ws.Range("K134").Value = 4363.5714
$ws.Range("M134").Value = -1828.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 243.85715
$ws.Range("I7").Value = 177
$ws.Range("K7").Value = 177
$ws.Range("M7").Value = -64

$ws.Range("H31").Value = 493552.53
$ws.Range("I31").Value = 10387.143
$ws.Range("K31").Value = 10387.143
$ws.Range("M31").Value = -10092.143

$ws.Range("H34").Value = 493552.53
$ws.Range("I34").Value = 10387.143
$ws.Range("K34").Value = 10387.143
$ws.Range("M34").Value = -10185.143

$ws.Range("H39").Value = 14123
$ws.Range("I39").Value = 14123
$ws.Range("K39").Value = 14123
$ws.Range("M39").Value = -13732

$ws.Range("H49").Value = 14123
$ws.Range("I49").Value = 14123
$ws.Range("K49").Value = 14123
$ws.Range("M49").Value = -13941

$ws.Range("H99").Value = 3036.3462
$ws.Range("J99").Value = 3114.25
$ws.Range("L99").Value = 3114.25
$ws.Range("N99").Value = -6110.25

$ws.Range("H126").Value = 3036.3462
$ws.Range("J126").Value = 3114.25
$ws.Range("L126").Value = 9342.75
$ws.Range("N126").Value = -14282.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 849.8182
$ws.Range("J104").Value = 999.6
$ws.Range("L104").Value = 2998.8
$ws.Range("N104").Value = -8240.799999999999

$ws.Range("H121").Value = 1946.1
$ws.Range("I121").Value = 482.5
$ws.Range("J121").Value = 2921.8333
$ws.Range("K121").Value = 1447.5
$ws.Range("L121").Value = 8765.499899999999
$ws.Range("M121").Value = -137.5
$ws.Range("N121").Value = -11385.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 432.11765
$ws.Range("I2").Value = 166.14285
$ws.Range("J2").Value = 618.3
$ws.Range("K2").Value = 166.14285
$ws.Range("L2").Value = 618.3
$ws.Range("M2").Value = -53.14285000000001
$ws.Range("N2").Value = -844.3

$ws.Range("H11").Value = 2985828.5
$ws.Range("J11").Value = 4121105.8
$ws.Range("L11").Value = 4121105.8
$ws.Range("N11").Value = -4121383.8

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H126").Value = 4624.4287
$ws.Range("I126").Value = 4377.8
$ws.Range("K126").Value = 13133.4
$ws.Range("M126").Value = -10663.4

$ws.Range("H132").Value = 27783316
$ws.Range("I132").Value = 33335542
$ws.Range("J132").Value = 22190.5
$ws.Range("K132").Value = 100006626
$ws.Range("L132").Value = 66571.5
$ws.Range("M132").Value = -100004096
$ws.Range("N132").Value = -71631.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 46921.043
$ws.Range("I7").Value = 2505.5
$ws.Range("K7").Value = 2505.5
$ws.Range("M7").Value = -2393.5

$ws.Range("H42").Value = 15000
$ws.Range("J42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -16126

$ws.Range("H46").Value = 4048.4194
$ws.Range("I46").Value = 2050.2778
$ws.Range("J46").Value = 6815.077
$ws.Range("K46").Value = 2050.2778
$ws.Range("L46").Value = 6815.077
$ws.Range("M46").Value = -1862.2778
$ws.Range("N46").Value = -7191.077

$ws.Range("H49").Value = 15000
$ws.Range("J49").Value = 15000
$ws.Range("L49").Value = 15000
$ws.Range("N49").Value = -15294

$ws.Range("H61").Value = 1238.4667
$ws.Range("I61").Value = 1238.4667
$ws.Range("K61").Value = 1238.4667
$ws.Range("M61").Value = -1036.4667

$ws.Range("H113").Value = 1238.4667
$ws.Range("I113").Value = 1238.4667
$ws.Range("K113").Value = 1238.4667
$ws.Range("M113").Value = 931.5333000000001

$ws.Range("H126").Value = 46921.043
$ws.Range("I126").Value = 2505.5
$ws.Range("K126").Value = 7516.5
$ws.Range("M126").Value = -5046.5

$ws.Range("H132").Value = 336262
$ws.Range("I132").Value = 2695.238
$ws.Range("K132").Value = 8085.714
$ws.Range("M132").Value = -5555.714

$ws.Range("H136").Value = 86120.53
$ws.Range("I136").Value = 11179.8
$ws.Range("K136").Value = 33539.39999999999
$ws.Range("M136").Value = -30989.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 19060028
$ws.Range("J75").Value = 18825034
$ws.Range("L75").Value = 18825034
$ws.Range("N75").Value = -18826906

$ws.Range("H78").Value = 19060028
$ws.Range("J78").Value = 18825034
$ws.Range("L78").Value = 56475102
$ws.Range("N78").Value = -56484462

$ws.Range("H132").Value = 3005.182
$ws.Range("I132").Value = 2506
$ws.Range("J132").Value = 6166.6665
$ws.Range("K132").Value = 7518
$ws.Range("L132").Value = 18499.9995
$ws.Range("M132").Value = -4988
$ws.Range("N132").Value = -23559.9995
